# Auto-generated edit script: updates crypto price/volume table cells to match the new snapshot.
# Numeric-looking "Price" values are written with a leading apostrophe so Excel stores them
# as text (preserving trailing zeros / exact formatting) instead of coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.867.78'
$ws.Range("E2").Value = '  -2.38%  '
$ws.Range("D3").Value = '3.471.75'
$ws.Range("E3").Value = '  +0.85%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''581.78'
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("D6").Value = '''171.73'
$ws.Range("E6").Value = '  -3.79%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '''0.596'
$ws.Range("E8").Value = '  -1.29%  '
$ws.Range("D9").Value = '3.473.16'
$ws.Range("E9").Value = '  +1.06%  '
$ws.Range("D10").Value = '''0.130'
$ws.Range("E10").Value = '  -5.72%  '
$ws.Range("E11").Value = '  -1.48%  '
$ws.Range("D12").Value = '''0.409'
$ws.Range("E12").Value = '  -3.84%  '
$ws.Range("D13").Value = '4.074.41'
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("E14").Value = '  +1.13%  '
$ws.Range("D15").Value = '''29.79'
$ws.Range("E15").Value = '  -6.44%  '
$ws.Range("D16").Value = '65.970.13'
$ws.Range("E16").Value = '  -2.19%  '
$ws.Range("D17").Value = '''0.0000170'
$ws.Range("E17").Value = '  -3.15%  '
$ws.Range("D18").Value = '3.478.14'
$ws.Range("E18").Value = '  +1.25%  '
$ws.Range("D19").Value = '''5.91'
$ws.Range("E19").Value = '  -3.24%  '
$ws.Range("D20").Value = '''13.85'
$ws.Range("E20").Value = '  -0.64%  '
$ws.Range("D21").Value = '''365.83'
$ws.Range("E21").Value = '  -4.97%  '
$ws.Range("D22").Value = '''7.70'
$ws.Range("E22").Value = '  -1.16%  '
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '''72.52'
$ws.Range("E23").Value = '  +2.07%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '''0.997'
$ws.Range("E24").Value = '  -0.25%  '
$ws.Range("E25").Value = '  +5.86%  '
$ws.Range("D26").Value = '''0.532'
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("D27").Value = '''9.56'
$ws.Range("E27").Value = '  -5.94%  '
$ws.Range("D28").Value = '''0.179'
$ws.Range("E28").Value = '  +2.46%  '
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").Value = '''24.00'
$ws.Range("E30").Value = '  +3.34%  '
$ws.Range("D31").Value = '''5.75'
$ws.Range("E31").Value = '  -4.79%  '
$ws.Range("D32").Value = '''1.98'
$ws.Range("E32").Value = '  -2.74%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").Value = '''7.11'
$ws.Range("E34").Value = '  -1.10%  '
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").Value = '''1.29'
$ws.Range("E35").Value = '  -6.18%  '
$ws.Range("E36").Value = '  -1.24%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '''159.78'
$ws.Range("E37").Value = '  -0.78%  '
$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").Value = '''29.28'
$ws.Range("E38").Value = '  +13.96%  '
$ws.Range("D39").Value = '''0.888'
$ws.Range("E39").Value = '  +0.98%  '
$ws.Range("D40").Value = '2.825.71'
$ws.Range("E40").Value = '  +4.83%  '
$ws.Range("D41").Value = '''1.75'
$ws.Range("E41").Value = '  -5.23%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '''6.42'
$ws.Range("E42").Value = '  -2.77%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '''4.43'
$ws.Range("E43").Value = '  -1.31%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = '''2.54'
$ws.Range("E44").Value = '  -6.92%  '
$ws.Range("D45").Value = '''0.0680'
$ws.Range("E45").Value = '  -4.12%  '
$ws.Range("D46").Value = '''40.09'
$ws.Range("E46").Value = '  -2.56%  '
$ws.Range("D47").Value = '''24.02'
$ws.Range("E47").Value = '  -7.00%  '
$ws.Range("D48").Value = '''0.0287'
$ws.Range("E48").Value = '  -2.74%  '
$ws.Range("D49").Value = '''324.80'
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").Value = '''0.811'
$ws.Range("E50").Value = '  -2.58%  '
$ws.Range("E51").Value = '  -2.67%  '
